$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (A1:Z1), columns AA1:AC1 removed ---
$headers = @("N","Tobs","d","variant","cauchy_pcombine","cauchy_cond_pcombine","bonferroni_pcombine","bonferroni_cond_pcombine","iu_pcombine","iu_cond_pcombine","Genmean_rneg_pcombine","Genmean_rneg_cond_pcombine","Genmean_pcombine","Genmean_cond_pcombine","Geomean_pcombine","Geomean_cond_pcombine","cauchy_bonf_pcombine","cauchy_bonf_cond_pcombine","iu_bonf_pcombine","iu_bonf_cond_pcombine","Genmean_rneg_bonf_pcombine","Genmean_rneg_bonf_cond_pcombine","Genmean_bonf_pcombine","Genmean_bonf_cond_pcombine","Geomean_bonf_pcombine","Geomean_bonf_cond_pcombine")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Remove now-unused trailing header cells AA1:AC1 (columns 27-29)
$ws.Cells.Item(1, 27).Value = $null
$ws.Cells.Item(1, 28).Value = $null
$ws.Cells.Item(1, 29).Value = $null

# --- Row 2 data (A2:Z2), columns AA2:AC2 removed ---
$row2 = @(80, 20, 0, "overall_holds", 0.0625, 0.04, 0.07, 0.0375, 0.1625, 0.12, 0.0675, 0.0325, 0.005, 0, 0.0075, 0, 0.0375, 0.0175, 0.0925, 0.0575, 0.0325, 0.0175, 0.0325, 0.0175, 0.0325, 0.0175)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$ws.Cells.Item(2, 27).Value = $null
$ws.Cells.Item(2, 28).Value = $null
$ws.Cells.Item(2, 29).Value = $null

# --- Row 3 data (A3:Z3), new row ---
$row3 = @(80, 50, 0, "overall_holds", 0.045, 0.045, 0.045, 0.0325, 0.1525, 0.1175, 0.0425, 0.03, 0.0025, 0.0025, 0.0025, 0.0025, 0.0225, 0.015, 0.0625, 0.05, 0.0225, 0.015, 0.0225, 0.015, 0.0225, 0.015)

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
